$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 326, shifting rows 326:398 down to 327:399.
$ws.Rows("326:326").Insert()

# Populate the newly inserted row 326 with the new record
# (same Mercado/Region/Categoria/Variedad/Calidad as its neighbours,
# new Fecha + Volumen/Precio data).
$ws.Range("A326").Value = 9
$ws.Range("B326").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C326").Value = "Metropolitana"
$ws.Range("D326").Value = 45211
$ws.Range("E326").Value = 13
$ws.Range("F326").Value = 100112001
$ws.Range("G326").Value = "Berenjena"
$ws.Range("H326").Value = "Sin especificar"
$ws.Range("I326").Value = "Primera"
$ws.Range("J326").Value = 70
$ws.Range("K326").Value = 9000
$ws.Range("L326").Value = 10000
$ws.Range("M326").Value = 9500
$ws.Range("N326").Value = "`$/caja 50 unidades"
$ws.Range("O326").Value = "Región de Arica y Parinacota"
$ws.Range("P326").Value = 190
$ws.Range("Q326").Value = 50
$ws.Range("R326").Value = "Hortaliza"
